$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.918.76"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.870.06"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5084"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07179"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8906"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.883.71"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07495"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.225"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008505"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "26.968.30"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.014"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "2.115.38"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.386"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.781"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.079"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.692"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.705"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09158"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05048"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.978"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.214"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.23%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.515"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5629"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.072"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.612"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.534"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1483"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4776"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9991"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.556"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
